$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 4; $r -le 11; $r++) {
    $ws.Range("A$r").Value = "mollie_ABC123"
    $ws.Range("B$r").Value = 1445758
    $ws.Range("C$r").Value = "Ben Gortemaker"
    $ws.Range("D$r").Value = 1
}
